$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "25÷2=12, 1"
$t.Cell(1, 2).Range.Text = "77÷5=15, 2"
$t.Cell(1, 3).Range.Text = "43÷5=8, 3"
$t.Cell(1, 4).Range.Text = "95÷9=10, 5"
$t.Cell(1, 5).Range.Text = "56÷7=8, 0"
$t.Cell(5, 1).Range.Text = "89÷6=14, 5"
$t.Cell(5, 2).Range.Text = "15÷7=2, 1"
$t.Cell(5, 3).Range.Text = "72÷9=8, 0"
$t.Cell(5, 4).Range.Text = "59÷8=7, 3"
$t.Cell(5, 5).Range.Text = "84÷2=42, 0"
$t.Cell(9, 1).Range.Text = "63÷2=31, 1"
$t.Cell(9, 2).Range.Text = "54÷6=9, 0"
$t.Cell(9, 3).Range.Text = "96÷8=12, 0"
$t.Cell(9, 4).Range.Text = "49÷2=24, 1"
$t.Cell(9, 5).Range.Text = "95÷9=10, 5"
$t.Cell(13, 1).Range.Text = "59÷9=6, 5"
$t.Cell(13, 2).Range.Text = "23÷6=3, 5"
$t.Cell(13, 3).Range.Text = "85÷3=28, 1"
$t.Cell(13, 4).Range.Text = "68÷8=8, 4"
$t.Cell(13, 5).Range.Text = "23÷3=7, 2"
$t.Cell(17, 1).Range.Text = "56÷5=11, 1"
$t.Cell(17, 2).Range.Text = "74÷7=10, 4"
$t.Cell(17, 3).Range.Text = "43÷3=14, 1"
$t.Cell(17, 4).Range.Text = "32÷4=8, 0"
$t.Cell(17, 5).Range.Text = "42÷7=6, 0"

Write-Host "Done"
